$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 4 new rows for "EVEN NUMBERS / PRIME NUMBERS / ENUMERATION / MERGE SORT"
#        right before the second table's header (old row 9), pushing everything below down by 4.
$ws.Rows("9:12").Insert()

# Fill the newly inserted rows (9-12) with index + algorithm name, mirroring the
# A/B layout used by the first table (rows 3-8).
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "EVEN NUMBERS"
$ws.Range("A10").Value = 2
$ws.Range("B10").Value = "PRIME NUMBERS"
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "ENUMERATION"
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = "MERGE SORT"

# --- 2. Remove the underline styling that used to sit on the second table's
#        index/name columns (old rows 10-15, now shifted to 14-19).
$ws.Range("A14:B19").Font.Underline = -4142

# --- 3. Append 4 more rows (20-23) to the bottom of the second table with the
#        same new algorithms, re-using the blank placeholder rows already there.
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "EVEN NUMBERS"
$ws.Range("A21").Value = 2
$ws.Range("B21").Value = "PRIME NUMBERS"
$ws.Range("A22").Value = 3
$ws.Range("B22").Value = "ENUMERATION"
$ws.Range("A23").Value = 4
$ws.Range("B23").Value = "MERGE SORT"

# --- 4. Restore the selection to the cell the author had active when saving.
$ws.Range("C8").Select() | Out-Null
